$d = $word.ActiveDocument

# --- Step 1: locate the exact end of "...normes de programmation." which is
# immediately followed by the (hidden) _GoBack bookmark and end of paragraph.
$fullText = $d.Content.Text
$anchor = "J’ai aussi définit nos normes de programmation."
$pos = $fullText.IndexOf($anchor)
$insertPoint = $pos + $anchor.Length

$r = $d.Range($insertPoint, $insertPoint)

# --- Step 2: insert the two new paragraphs' text as one run of plain text,
# using a placeholder marker where the paragraph breaks should go (so the
# "_GoBack" bookmark, which sits right at $insertPoint, ends up after all of
# the newly-typed text instead of swallowing the later paragraphs).
$marker = "||PBREAK||"
$newText = $marker + "Le 18 novembre 2014" + $marker + "PLACEHOLDER_TOOLBAR_PARA"
$r.InsertBefore($newText)

# Turn the markers into real paragraph marks via Find/Replace (keeps the
# bookmark anchored after the final character of inserted text).
$d.Content.Find.Execute($marker, $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2)

# --- Step 3: replace the placeholder paragraph (now the last paragraph in
# the document) with the fully-formed run/proofErr structure, including the
# _GoBack bookmark markers at the end, matching the target markup exactly.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $lastPara.Range

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">J’ai refait notre diagramme de classe au complet suite à un problème où tout a été supprimé et aucun retour n’était possible, il est plus à jour. J’ai aussi entamé la classe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>toolbar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pour afficher les outils disponibles.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)

# --- Step 4: InsertXML on a whole last-paragraph Range leaves a stray empty
# paragraph behind (mirrors InsertParagraphAfter semantics at story end).
# Remove it by deleting the paragraph mark that now separates our finished
# paragraph from that trailing empty one.
$trailingText = $d.Paragraphs($d.Paragraphs.Count).Range.Text
if ($trailingText.Length -le 1) {
    $prevPara = $d.Paragraphs($d.Paragraphs.Count - 1)
    $markRange = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
    $markRange.Delete()
}

Write-Output "done"
